# Applies weekly-refresh data update to Hortaliza, Vega Modelo de Temuco - Camote sheet
# (columns Fecha, Volumen, Precio minimo/maximo/promedio, Unidad, Origen, Precio $/Kg, Kg o Unidades)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44340
$ws.Cells.Item(2, 10).Value = 40
$ws.Cells.Item(2, 11).Value = 18000
$ws.Cells.Item(2, 12).Value = 18000
$ws.Cells.Item(2, 13).Value = 18000
$ws.Cells.Item(2, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(2, 15).Value = 'Perú'
$ws.Cells.Item(2, 16).Value = 900
$ws.Cells.Item(2, 17).Value = 20

# Row 3
$ws.Cells.Item(3, 4).Value = 44389
$ws.Cells.Item(3, 10).Value = 45
$ws.Cells.Item(3, 11).Value = 20000
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = 20000
$ws.Cells.Item(3, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(3, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(3, 16).Value = 1333
$ws.Cells.Item(3, 17).Value = 15

# Row 4
$ws.Cells.Item(4, 4).Value = 44188
$ws.Cells.Item(4, 10).Value = 20
$ws.Cells.Item(4, 11).Value = 20000
$ws.Cells.Item(4, 12).Value = 20000
$ws.Cells.Item(4, 13).Value = 20000
$ws.Cells.Item(4, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(4, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(4, 16).Value = 1333
$ws.Cells.Item(4, 17).Value = 15

# Row 5
$ws.Cells.Item(5, 4).Value = 44424
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(5, 11).Value = 20000
$ws.Cells.Item(5, 12).Value = 20000
$ws.Cells.Item(5, 13).Value = 20000
$ws.Cells.Item(5, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(5, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(5, 16).Value = 1333
$ws.Cells.Item(5, 17).Value = 15

# Row 6
$ws.Cells.Item(6, 4).Value = 44369
$ws.Cells.Item(6, 10).Value = 20
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 20000
$ws.Cells.Item(6, 13).Value = 20000
$ws.Cells.Item(6, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(6, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(6, 16).Value = 1333
$ws.Cells.Item(6, 17).Value = 15

# Row 7
$ws.Cells.Item(7, 4).Value = 44369
$ws.Cells.Item(7, 10).Value = 20
$ws.Cells.Item(7, 11).Value = 20000
$ws.Cells.Item(7, 12).Value = 20000
$ws.Cells.Item(7, 13).Value = 20000
$ws.Cells.Item(7, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(7, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(7, 16).Value = 1000
$ws.Cells.Item(7, 17).Value = 20

# Row 8
$ws.Cells.Item(8, 4).Value = 44341
$ws.Cells.Item(8, 10).Value = 40
$ws.Cells.Item(8, 11).Value = 17000
$ws.Cells.Item(8, 12).Value = 18000
$ws.Cells.Item(8, 13).Value = 17500
$ws.Cells.Item(8, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(8, 15).Value = 'Perú'
$ws.Cells.Item(8, 16).Value = 875
$ws.Cells.Item(8, 17).Value = 20

# Row 9
$ws.Cells.Item(9, 4).Value = 44525
$ws.Cells.Item(9, 10).Value = 40
$ws.Cells.Item(9, 11).Value = 20000
$ws.Cells.Item(9, 12).Value = 20000
$ws.Cells.Item(9, 13).Value = 20000
$ws.Cells.Item(9, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(9, 15).Value = 'Perú'
$ws.Cells.Item(9, 16).Value = 1333
$ws.Cells.Item(9, 17).Value = 15

# Row 10
$ws.Cells.Item(10, 4).Value = 44448
$ws.Cells.Item(10, 10).Value = 45
$ws.Cells.Item(10, 11).Value = 20000
$ws.Cells.Item(10, 12).Value = 20000
$ws.Cells.Item(10, 13).Value = 20000
$ws.Cells.Item(10, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(10, 15).Value = 'Perú'
$ws.Cells.Item(10, 16).Value = 1000
$ws.Cells.Item(10, 17).Value = 20

# Row 11
$ws.Cells.Item(11, 4).Value = 44294
$ws.Cells.Item(11, 10).Value = 5
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 13).Value = 20000
$ws.Cells.Item(11, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(11, 15).Value = 'Perú'
$ws.Cells.Item(11, 16).Value = 1333
$ws.Cells.Item(11, 17).Value = 15

# Row 12
$ws.Cells.Item(12, 4).Value = 44497
$ws.Cells.Item(12, 10).Value = 30
$ws.Cells.Item(12, 11).Value = 20000
$ws.Cells.Item(12, 12).Value = 20000
$ws.Cells.Item(12, 13).Value = 20000
$ws.Cells.Item(12, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(12, 15).Value = 'Perú'
$ws.Cells.Item(12, 16).Value = 1333
$ws.Cells.Item(12, 17).Value = 15

# Row 13
$ws.Cells.Item(13, 4).Value = 44497
$ws.Cells.Item(13, 10).Value = 40
$ws.Cells.Item(13, 11).Value = 20000
$ws.Cells.Item(13, 12).Value = 20000
$ws.Cells.Item(13, 13).Value = 20000
$ws.Cells.Item(13, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(13, 15).Value = 'Perú'
$ws.Cells.Item(13, 16).Value = 1000
$ws.Cells.Item(13, 17).Value = 20

# Row 14
$ws.Cells.Item(14, 4).Value = 44179
$ws.Cells.Item(14, 10).Value = 20
$ws.Cells.Item(14, 11).Value = 20000
$ws.Cells.Item(14, 12).Value = 20000
$ws.Cells.Item(14, 13).Value = 20000
$ws.Cells.Item(14, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(14, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(14, 16).Value = 1333
$ws.Cells.Item(14, 17).Value = 15

# Row 15
$ws.Cells.Item(15, 4).Value = 44425
$ws.Cells.Item(15, 10).Value = 10
$ws.Cells.Item(15, 11).Value = 20000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 13).Value = 20000
$ws.Cells.Item(15, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(15, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(15, 16).Value = 1333
$ws.Cells.Item(15, 17).Value = 15

# Row 16
$ws.Cells.Item(16, 4).Value = 44498
$ws.Cells.Item(16, 10).Value = 20
$ws.Cells.Item(16, 11).Value = 20000
$ws.Cells.Item(16, 12).Value = 20000
$ws.Cells.Item(16, 13).Value = 20000
$ws.Cells.Item(16, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(16, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(16, 16).Value = 1000
$ws.Cells.Item(16, 17).Value = 20

# Row 17
$ws.Cells.Item(17, 4).Value = 44579
$ws.Cells.Item(17, 10).Value = 50
$ws.Cells.Item(17, 11).Value = 20000
$ws.Cells.Item(17, 12).Value = 20000
$ws.Cells.Item(17, 13).Value = 20000
$ws.Cells.Item(17, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(17, 15).Value = 'Perú'
$ws.Cells.Item(17, 16).Value = 1000
$ws.Cells.Item(17, 17).Value = 20

# Row 18
$ws.Cells.Item(18, 4).Value = 44578
$ws.Cells.Item(18, 10).Value = 50
$ws.Cells.Item(18, 11).Value = 20000
$ws.Cells.Item(18, 12).Value = 20000
$ws.Cells.Item(18, 13).Value = 20000
$ws.Cells.Item(18, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(18, 15).Value = 'Perú'
$ws.Cells.Item(18, 16).Value = 1000
$ws.Cells.Item(18, 17).Value = 20

# Row 19
$ws.Cells.Item(19, 4).Value = 44329
$ws.Cells.Item(19, 10).Value = 40
$ws.Cells.Item(19, 11).Value = 20000
$ws.Cells.Item(19, 12).Value = 20000
$ws.Cells.Item(19, 13).Value = 20000
$ws.Cells.Item(19, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(19, 15).Value = 'Perú'
$ws.Cells.Item(19, 16).Value = 1333
$ws.Cells.Item(19, 17).Value = 15

# Row 20
$ws.Cells.Item(20, 4).Value = 44508
$ws.Cells.Item(20, 10).Value = 40
$ws.Cells.Item(20, 11).Value = 20000
$ws.Cells.Item(20, 12).Value = 20000
$ws.Cells.Item(20, 13).Value = 20000
$ws.Cells.Item(20, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(20, 15).Value = 'Perú'
$ws.Cells.Item(20, 16).Value = 1333
$ws.Cells.Item(20, 17).Value = 15

# Row 21
$ws.Cells.Item(21, 4).Value = 44455
$ws.Cells.Item(21, 10).Value = 30
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 20000
$ws.Cells.Item(21, 13).Value = 20000
$ws.Cells.Item(21, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(21, 15).Value = 'Perú'
$ws.Cells.Item(21, 16).Value = 1000
$ws.Cells.Item(21, 17).Value = 20

# Row 22
$ws.Cells.Item(22, 4).Value = 44466
$ws.Cells.Item(22, 10).Value = 20
$ws.Cells.Item(22, 11).Value = 25000
$ws.Cells.Item(22, 12).Value = 25000
$ws.Cells.Item(22, 13).Value = 25000
$ws.Cells.Item(22, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(22, 15).Value = 'Perú'
$ws.Cells.Item(22, 16).Value = 1667
$ws.Cells.Item(22, 17).Value = 15

# Row 23
$ws.Cells.Item(23, 4).Value = 44364
$ws.Cells.Item(23, 10).Value = 15
$ws.Cells.Item(23, 11).Value = 20000
$ws.Cells.Item(23, 12).Value = 20000
$ws.Cells.Item(23, 13).Value = 20000
$ws.Cells.Item(23, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(23, 15).Value = 'Perú'
$ws.Cells.Item(23, 16).Value = 1333
$ws.Cells.Item(23, 17).Value = 15

# Row 24
$ws.Cells.Item(24, 4).Value = 44452
$ws.Cells.Item(24, 10).Value = 50
$ws.Cells.Item(24, 11).Value = 20000
$ws.Cells.Item(24, 12).Value = 20000
$ws.Cells.Item(24, 13).Value = 20000
$ws.Cells.Item(24, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(24, 15).Value = 'Perú'
$ws.Cells.Item(24, 16).Value = 1000
$ws.Cells.Item(24, 17).Value = 20

# Row 25
$ws.Cells.Item(25, 4).Value = 44512
$ws.Cells.Item(25, 10).Value = 30
$ws.Cells.Item(25, 11).Value = 20000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 13).Value = 20000
$ws.Cells.Item(25, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(25, 15).Value = 'Perú'
$ws.Cells.Item(25, 16).Value = 1000
$ws.Cells.Item(25, 17).Value = 20

# Row 26
$ws.Cells.Item(26, 4).Value = 44511
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 20000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = 20000
$ws.Cells.Item(26, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(26, 15).Value = 'Perú'
$ws.Cells.Item(26, 16).Value = 1000
$ws.Cells.Item(26, 17).Value = 20

# Row 27
$ws.Cells.Item(27, 4).Value = 44175
$ws.Cells.Item(27, 10).Value = 20
$ws.Cells.Item(27, 11).Value = 20000
$ws.Cells.Item(27, 12).Value = 20000
$ws.Cells.Item(27, 13).Value = 20000
$ws.Cells.Item(27, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(27, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(27, 16).Value = 1333
$ws.Cells.Item(27, 17).Value = 15

# Row 28
$ws.Cells.Item(28, 4).Value = 44161
$ws.Cells.Item(28, 10).Value = 20
$ws.Cells.Item(28, 11).Value = 20000
$ws.Cells.Item(28, 12).Value = 20000
$ws.Cells.Item(28, 13).Value = 20000
$ws.Cells.Item(28, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(28, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(28, 16).Value = 1333
$ws.Cells.Item(28, 17).Value = 15

# Row 29
$ws.Cells.Item(29, 4).Value = 44567
$ws.Cells.Item(29, 10).Value = 50
$ws.Cells.Item(29, 11).Value = 18000
$ws.Cells.Item(29, 12).Value = 18000
$ws.Cells.Item(29, 13).Value = 18000
$ws.Cells.Item(29, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(29, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(29, 16).Value = 900
$ws.Cells.Item(29, 17).Value = 20

# Row 30
$ws.Cells.Item(30, 4).Value = 44532
$ws.Cells.Item(30, 10).Value = 40
$ws.Cells.Item(30, 11).Value = 18000
$ws.Cells.Item(30, 12).Value = 18000
$ws.Cells.Item(30, 13).Value = 18000
$ws.Cells.Item(30, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(30, 15).Value = 'Perú'
$ws.Cells.Item(30, 16).Value = 900
$ws.Cells.Item(30, 17).Value = 20

# Row 31
$ws.Cells.Item(31, 4).Value = 44441
$ws.Cells.Item(31, 10).Value = 40
$ws.Cells.Item(31, 11).Value = 20000
$ws.Cells.Item(31, 12).Value = 20000
$ws.Cells.Item(31, 13).Value = 20000
$ws.Cells.Item(31, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(31, 15).Value = 'Perú'
$ws.Cells.Item(31, 16).Value = 1000
$ws.Cells.Item(31, 17).Value = 20

# Row 32
$ws.Cells.Item(32, 4).Value = 44316
$ws.Cells.Item(32, 10).Value = 20
$ws.Cells.Item(32, 11).Value = 20000
$ws.Cells.Item(32, 12).Value = 20000
$ws.Cells.Item(32, 13).Value = 20000
$ws.Cells.Item(32, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(32, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(32, 16).Value = 1333
$ws.Cells.Item(32, 17).Value = 15

# Row 33
$ws.Cells.Item(33, 4).Value = 44315
$ws.Cells.Item(33, 10).Value = 30
$ws.Cells.Item(33, 11).Value = 20000
$ws.Cells.Item(33, 12).Value = 20000
$ws.Cells.Item(33, 13).Value = 20000
$ws.Cells.Item(33, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(33, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(33, 16).Value = 1333
$ws.Cells.Item(33, 17).Value = 15

# Row 34
$ws.Cells.Item(34, 4).Value = 44315
$ws.Cells.Item(34, 10).Value = 30
$ws.Cells.Item(34, 11).Value = 20000
$ws.Cells.Item(34, 12).Value = 20000
$ws.Cells.Item(34, 13).Value = 20000
$ws.Cells.Item(34, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(34, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(34, 16).Value = 1000
$ws.Cells.Item(34, 17).Value = 20

# Row 35
$ws.Cells.Item(35, 4).Value = 44186
$ws.Cells.Item(35, 10).Value = 20
$ws.Cells.Item(35, 11).Value = 20000
$ws.Cells.Item(35, 12).Value = 20000
$ws.Cells.Item(35, 13).Value = 20000
$ws.Cells.Item(35, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(35, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(35, 16).Value = 1333
$ws.Cells.Item(35, 17).Value = 15

# Row 36
$ws.Cells.Item(36, 4).Value = 44496
$ws.Cells.Item(36, 10).Value = 30
$ws.Cells.Item(36, 11).Value = 20000
$ws.Cells.Item(36, 12).Value = 20000
$ws.Cells.Item(36, 13).Value = 20000
$ws.Cells.Item(36, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(36, 15).Value = 'Perú'
$ws.Cells.Item(36, 16).Value = 1000
$ws.Cells.Item(36, 17).Value = 20

# Row 37
$ws.Cells.Item(37, 4).Value = 44438
$ws.Cells.Item(37, 10).Value = 40
$ws.Cells.Item(37, 11).Value = 20000
$ws.Cells.Item(37, 12).Value = 20000
$ws.Cells.Item(37, 13).Value = 20000
$ws.Cells.Item(37, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(37, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(37, 16).Value = 1333
$ws.Cells.Item(37, 17).Value = 15

# Row 38
$ws.Cells.Item(38, 4).Value = 44385
$ws.Cells.Item(38, 10).Value = 18
$ws.Cells.Item(38, 11).Value = 20000
$ws.Cells.Item(38, 12).Value = 20000
$ws.Cells.Item(38, 13).Value = 20000
$ws.Cells.Item(38, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(38, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(38, 16).Value = 1000
$ws.Cells.Item(38, 17).Value = 20

# Row 39
$ws.Cells.Item(39, 4).Value = 44529
$ws.Cells.Item(39, 10).Value = 15
$ws.Cells.Item(39, 11).Value = 20000
$ws.Cells.Item(39, 12).Value = 20000
$ws.Cells.Item(39, 13).Value = 20000
$ws.Cells.Item(39, 14).Value = '$/malla 20 kilos'
$ws.Cells.Item(39, 15).Value = 'Perú'
$ws.Cells.Item(39, 16).Value = 1000
$ws.Cells.Item(39, 17).Value = 20

# Row 40
$ws.Cells.Item(40, 4).Value = 44321
$ws.Cells.Item(40, 10).Value = 15
$ws.Cells.Item(40, 11).Value = 25000
$ws.Cells.Item(40, 12).Value = 25000
$ws.Cells.Item(40, 13).Value = 25000
$ws.Cells.Item(40, 14).Value = '$/caja 15 kilos granel'
$ws.Cells.Item(40, 15).Value = 'Perú'
$ws.Cells.Item(40, 16).Value = 1667
$ws.Cells.Item(40, 17).Value = 15

Write-Host "Applied weekly data refresh"
